# Coding Test Excel 1 - fill in answers for tasks 2 and 3, and clean up the
# scratch VLOOKUP helper table that was used to work out task 3's answer.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Task 2: average sales amount for the South region.
$ws2.Range("F2").Formula = '=AVERAGEIF(Sheet1!F1:F11,"South",Sheet1!G1:G11)'

# Task 3: maximum sales amount for the West region.
$ws2.Range("F3").Formula = '=MAXIFS(Sheet1!G1:G11,Sheet1!F1:F11,"West")'

# Remove the scratch helper table (rows 6-7, columns A-E) that was used to
# manually look up the West-region values, keeping the existing cell
# formatting intact.
$ws2.Range("A6:E7").ClearContents()

# After typing the answer into F3 the cursor naturally lands on F4.
$ws2.Range("F4").Select()
